# "Add files via upload" — adds a new "Housing Material" comparison table to
# the "Part Comparison" sheet (rows 18 and 70-86), tweaks a couple of column
# widths, updates the F2 input (and dependent formulas) on the
# "Battery Life Estimate" sheet, and moves the active sheet/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Part Comparison
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Part Comparison")

# New two-cell row appended to the "Air Pumps" table
$ws.Range("A18").Value = "Pump Type"
$ws.Range("B18").Value = "Pressure"

# New "Housing Material" comparison table ----------------------------------

# Section title (copy format from another section title, e.g. A22 "Valves")
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A70").PasteSpecial(-4122) | Out-Null
$ws.Range("A70").Value = "Housing Material"

# Header row (copy format from another header row, e.g. A23:D23)
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A71:H71").PasteSpecial(-4122) | Out-Null
$ws.Range("A71").Value = "(1-10 scale)"
$ws.Range("B71").Value = "ABS"
$ws.Range("C71").Value = "Acrylic"
$ws.Range("D71").Value = "PP"
$ws.Range("E71").Value = "Stainless Steel"
$ws.Range("F71").Value = "PETG"
$ws.Range("G71").Value = "Aluminum"
$ws.Range("H71").Value = "Titanium"

# Ranking row
$ws.Range("A72").Value = "Ranking"
$ws.Range("E72").Value = 2
$ws.Range("F72").Value = 1

# Property rows
$ws.Range("A73").Value = "Environmentally Friendly"
$ws.Range("B73").Value = "No"

$ws.Range("A74").Value = "Smell"
$ws.Range("B74").Value = "Bad"

$ws.Range("A75").Value = "Fatigue Resistance"
$ws.Range("D75").Value = "High"

$ws.Range("A76").Value = "Impact Resistance"
$ws.Range("B76").Value = "High"

$ws.Range("A77").Value = "Durability"
$ws.Range("B77").Value = "High"
$ws.Range("D77").Value = "High"

$ws.Range("A78").Value = "UV Resistance"
$ws.Range("B78").Value = "Moderate (becomes brittle)"
$ws.Range("B78").WrapText = $true
$ws.Rows.Item(78).RowHeight = 45

$ws.Range("A79").Value = "Appearance"
$ws.Range("B79").Value = "Dull/matte"

$ws.Range("A80").Value = "Manufacturability"
$ws.Range("B80").Value = "Difficult"

$ws.Range("A81").Value = "Biocompatablity"

$ws.Range("A82").Value = "Hypoallergenic"
$ws.Range("B82").Value = "Low/Moderate"

$ws.Range("A83").Value = "Hygiene"
$ws.Range("D83").Value = "Excellent"

$ws.Range("A84").Value = "Paint/Finish Adhesion"
$ws.Range("D84").Value = "Poor"

$ws.Range("A85").Value = "Biocompatablity"
$ws.Range("C36").Copy() | Out-Null
$ws.Range("B85:C85").PasteSpecial(-4122) | Out-Null
$ws.Range("D85").Value = "Excellent"

$ws.Range("A86").Value = "Heat Resistance"
$ws.Range("B86").Value = "High"
$ws.Range("D86").Value = "Moderate"

# Column width tweaks
$ws.Columns.Item(1).ColumnWidth = 22.8
$ws.Columns.Item(2).ColumnWidth = 13.15
$ws.Columns.Item(5).ColumnWidth = 13.15

# ---------------------------------------------------------------------------
# Sheet: Battery Life Estimate
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Battery Life Estimate")

$ws3.Rows.Item(1).RowHeight = 45
$ws3.Range("F2").Value = 800

$ws3.Range("H11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Make "Part Comparison" the active sheet/tab with its own selection
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("H76").Select() | Out-Null
